# Insert a new data row at sheet row 886 (pushes the existing rows 886-927
# down to 887-928) and populate it with the new sample:
#   2026/02/26, 木 (Thursday), hour 6, ranking 181
#
# This mirrors the upstream diff, which shows every row from 886 through
# the former last row (927) shifting down by one, with the brand-new row
# carrying the 2026/02/26 / 木 / 6 / 181 values and the sheet's used range
# growing from A1:D927 to A1:D928.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 886..927 down one position, leaving a blank row 886 to fill in.
$ws.Range("A886").EntireRow.Insert()

# Column A holds dates formatted as plain text (e.g. "2026/02/26") in this
# sheet, not real date serials. Force text number formatting before writing
# the value so Excel doesn't auto-convert the "YYYY/MM/DD"-looking string
# into a date serial, then clear the formatting again afterwards so the new
# row matches its neighbours (which carry no explicit cell style).
$dateCell = $ws.Cells.Item(886, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026/02/26"
$dateCell.ClearFormats()

$ws.Range("B886").Value = "木"
$ws.Range("C886").Value = 6
$ws.Range("D886").Value = 181
